# Remove the "MARIANA" balance row (account 004525587) from the Export
# sheet. This is row 2 (directly under the header row). Deleting the
# entire row shifts every following row up by one, matching the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

$ws.Rows.Item(2).Delete()
